$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "La grandezza dei puntini..."
# sentence (the one that also still carries the trailing "." run and the
# _GoBack bookmark).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*La grandezza dei puntini*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph 'La grandezza dei puntini...'"
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replace the single paragraph with four paragraphs:
#  1) the original sentence, now ending with a period inside the same run
#  2) "Quando un puntino viene cancellato, il numero di ogni puntino viene
#     aggiornato e non è più"
#  3) "come prima."
#  4) "Prima quando si selezionava un puntino un po' più grande dovevi per
#     forza prendere la parte centrale del puntino, ora è sistemato." with
#     the _GoBack bookmark (and trailing "." run) moved here, exactly as
#     it used to sit at the end of the original paragraph.
$xml = @"
<w:p $wNs w:rsidR="00FE0E2A" w:rsidRDefault="00C17CFF" w:rsidP="00AB580C"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>La grandezza dei puntini ora può cambiare anche dopo averli messi, basta selezionare un puntino e usare lo slider.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Quando un puntino viene cancellato, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>il numero di ogni puntino viene aggiornato e non è più</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>come prima.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Prima quando si selezionava un puntino un po’ più grande dovevi per forza prendere la parte centrale del puntino, ora è sistemato</w:t></w:r><w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>.</w:t></w:r></w:p>
"@

$target.Range.InsertXML($xml)

Write-Output "Paragraph split/updated successfully."
